$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the last row of names ("Krako"/"Krakic") with new ones ("Wynell"/"Aufderhar")
$ws.Range("A5").Value = "Wynell"
$ws.Range("B5").Value = "Aufderhar"

# Update the active selection to the next block where new names would be appended
$ws.Range("A6:B16").Select()
